$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Population size (P): 4190 -> 85000
$ws.Range("F3").Value = 85000

# Expected Occurrence (p): 0.5 -> 0.25
$ws.Range("F5").Value = 0.25

# Size of sample (n): 45 -> 42000
$ws.Range("E20").Value = 42000

# # matching criteria: 2 -> 13000
$ws.Range("E21").Value = 13000

# Update the active selection/view to match the final state (scrolled to E22)
$ws.Range("E22").Select()
